$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 124.44444
$ws.Range("I33").Value = 136.875
$ws.Range("K33").Value = 136.875
$ws.Range("M33").Value = 92.125
$ws.Range("H40").Value = 4808.4546
$ws.Range("I40").Value = 4253.643
$ws.Range("J40").Value = 5779.375
$ws.Range("K40").Value = 4253.643
$ws.Range("L40").Value = 5779.375
$ws.Range("M40").Value = -4078.643
$ws.Range("N40").Value = -6129.375
$ws.Range("H76").Value = 4699.6
$ws.Range("I76").Value = 3666
$ws.Range("K76").Value = 3666
$ws.Range("M76").Value = -3351
$ws.Range("H79").Value = 4699.6
$ws.Range("I79").Value = 3666
$ws.Range("K79").Value = 3666
$ws.Range("M79").Value = -2574
$ws.Range("H99").Value = 3315
$ws.Range("I99").Value = 1000
$ws.Range("J99").Value = 3893.75
$ws.Range("K99").Value = 3000
$ws.Range("L99").Value = 11681.25
$ws.Range("M99").Value = -1502
$ws.Range("N99").Value = -14677.25
$ws.Range("H116").Value = 4325
$ws.Range("J116").Value = 4325
$ws.Range("L116").Value = 4325
$ws.Range("N116").Value = -11209
$ws.Range("H135").Value = 1629.8462
$ws.Range("I135").Value = 1185.5
$ws.Range("J135").Value = 3111
$ws.Range("K135").Value = 10669.5
$ws.Range("L135").Value = 27999
$ws.Range("M135").Value = -8134.5
$ws.Range("N135").Value = -33069

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5828.1113
$ws.Range("I2").Value = 3494.5293
$ws.Range("K2").Value = 3494.5293
$ws.Range("M2").Value = -3381.5293
$ws.Range("H61").Value = 5000
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("H116").Value = 5828.1113
$ws.Range("I116").Value = 3494.5293
$ws.Range("K116").Value = 3494.5293
$ws.Range("M116").Value = -1200.5293
$ws.Range("H122").Value = 4048
$ws.Range("I122").Value = 3858.4
$ws.Range("K122").Value = 11575.2
$ws.Range("M122").Value = -9125.200000000001
$ws.Range("H136").Value = 5000
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5828.1113
$ws.Range("I3").Value = 3494.5293
$ws.Range("K3").Value = 3494.5293
$ws.Range("M3").Value = -3380.5293
$ws.Range("H86").Value = 3564.1052
$ws.Range("I86").Value = 1468.375
$ws.Range("K86").Value = 1468.375
$ws.Range("M86").Value = -345.375
$ws.Range("H89").Value = 3564.1052
$ws.Range("I89").Value = 1468.375
$ws.Range("K89").Value = 7341.875
$ws.Range("M89").Value = -1725.875
$ws.Range("H134").Value = 2775.9333
$ws.Range("I134").Value = 2369.8
$ws.Range("J134").Value = 3588.2
$ws.Range("K134").Value = 7109.400000000001
$ws.Range("L134").Value = 10764.6
$ws.Range("M134").Value = -4574.400000000001
$ws.Range("N134").Value = -15834.6
$ws.Range("H16").Value = 1449.5
$ws.Range("I16").Value = 1339.4
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 1339.4
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -1052.4
$ws.Range("N16").Value = -2574

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8920
$ws.Range("I31").Value = 3447.5
$ws.Range("K31").Value = 3447.5
$ws.Range("M31").Value = -3152.5
$ws.Range("H34").Value = 8920
$ws.Range("I34").Value = 3447.5
$ws.Range("K34").Value = 3447.5
$ws.Range("M34").Value = -3245.5
$ws.Range("H58").Value = 5095.9
$ws.Range("I58").Value = 998.6667
$ws.Range("K58").Value = 998.6667
$ws.Range("M58").Value = -795.6667
$ws.Range("H107").Value = 264.44446
$ws.Range("I107").Value = 268.13333
$ws.Range("K107").Value = 268.13333
$ws.Range("M107").Value = 1651.86667
$ws.Range("H113").Value = 1449.5
$ws.Range("I113").Value = 1339.4
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1339.4
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 830.5999999999999
$ws.Range("N113").Value = -6340
$ws.Range("H132").Value = 4091.138
$ws.Range("I132").Value = 3073.5
$ws.Range("K132").Value = 9220.5
$ws.Range("M132").Value = -6690.5
$ws.Range("H134").Value = 1895.75
$ws.Range("J134").Value = 1532.5
$ws.Range("L134").Value = 4597.5
$ws.Range("N134").Value = -9667.5
$ws.Range("H135").Value = 100000
$ws.Range("J135").Value = 100000
$ws.Range("L135").Value = 100000
$ws.Range("N135").Value = -110140
$ws.Range("H136").Value = 5095.9
$ws.Range("I136").Value = 998.6667
$ws.Range("K136").Value = 2996.0001
$ws.Range("M136").Value = -446.0001000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 358.8889
$ws.Range("I38").Value = 400.625
$ws.Range("K38").Value = 1201.875
$ws.Range("M38").Value = -854.875
$ws.Range("H132").Value = 1732.1111
$ws.Range("I132").Value = 995
$ws.Range("K132").Value = 8955
$ws.Range("M132").Value = -6425

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 32999
$ws.Range("I57").Value = 20000
$ws.Range("K57").Value = 20000
$ws.Range("M57").Value = -19180
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("N74").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("N77").Value = 0
$ws.Range("H126").Value = 1500
$ws.Range("I126").Value = 1500
$ws.Range("K126").Value = 4500
$ws.Range("M126").Value = -2030
$ws.Range("H140").Value = 135939.2
$ws.Range("J140").Value = 99999.75
$ws.Range("L140").Value = 99999.75
$ws.Range("N140").Value = -110359.75
$ws.Range("L74").ClearContents()
$ws.Range("L77").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6278.5835
$ws.Range("I7").Value = 4777.2856
$ws.Range("K7").Value = 4777.2856
$ws.Range("M7").Value = -4665.2856
$ws.Range("H16").Value = 712.1667
$ws.Range("I16").Value = 712.1667
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 712.1667
$ws.Range("L16").Value = 0
$ws.Range("N16").Value = -542.1667
$ws.Range("H22").Value = 843.1429000000001
$ws.Range("I22").Value = 780.4
$ws.Range("K22").Value = 780.4
$ws.Range("M22").Value = -485.4
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("N23").Value = 0
$ws.Range("H27").Value = 843.1429000000001
$ws.Range("I27").Value = 780.4
$ws.Range("K27").Value = 780.4
$ws.Range("M27").Value = -673.4
$ws.Range("H30").Value = 500
$ws.Range("I30").Value = 500
$ws.Range("K30").Value = 500
$ws.Range("M30").Value = -392
$ws.Range("H61").Value = 3779.4167
$ws.Range("I61").Value = 1669.125
$ws.Range("K61").Value = 1669.125
$ws.Range("M61").Value = -1467.125
$ws.Range("H113").Value = 3779.4167
$ws.Range("I113").Value = 1669.125
$ws.Range("K113").Value = 1669.125
$ws.Range("M113").Value = 500.875
$ws.Range("H126").Value = 6278.5835
$ws.Range("I126").Value = 4777.2856
$ws.Range("K126").Value = 14331.8568
$ws.Range("M126").Value = -11861.8568
$ws.Range("H136").Value = 3094
$ws.Range("I136").Value = 2377.75
$ws.Range("K136").Value = 7133.25
$ws.Range("M136").Value = -4583.25
$ws.Range("M16").ClearContents()
$ws.Range("L23").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 57500
$ws.Range("J64").Value = 57500
$ws.Range("L64").Value = 57500
$ws.Range("N64").Value = -57996
$ws.Range("H67").Value = 57500
$ws.Range("J67").Value = 57500
$ws.Range("L67").Value = 57500
$ws.Range("N67").Value = -59216
$ws.Range("H75").Value = 34000
$ws.Range("I75").Value = 34000
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 34000
$ws.Range("L75").Value = 0
$ws.Range("N75").Value = -33064
$ws.Range("H78").Value = 34000
$ws.Range("I78").Value = 34000
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 102000
$ws.Range("L78").Value = 0
$ws.Range("N78").Value = -97320
$ws.Range("H80").Value = 35000
$ws.Range("J80").Value = 35000
$ws.Range("L80").Value = 35000
$ws.Range("N80").Value = -36996
$ws.Range("H83").Value = 35000
$ws.Range("J83").Value = 35000
$ws.Range("L83").Value = 105000
$ws.Range("N83").Value = -114984
$ws.Range("H107").Value = 920
$ws.Range("I107").Value = 840
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 2520
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = -600
$ws.Range("N107").Value = -6840
$ws.Range("H126").Value = 5433.4443
$ws.Range("I126").Value = 2225.25
$ws.Range("K126").Value = 6675.75
$ws.Range("M126").Value = -4205.75
$ws.Range("H132").Value = 2854.0476
$ws.Range("I132").Value = 1513.5
$ws.Range("K132").Value = 4540.5
$ws.Range("M132").Value = -2010.5
$ws.Range("M75").ClearContents()
$ws.Range("M78").ClearContents()
